$wb = $excel.ActiveWorkbook

# --- Locate the existing "总计" (Total) sheet; it is the last sheet (rId5 / sheetId 5) ---
$total = $wb.Worksheets.Item($wb.Worksheets.Count)

# --- Step 1: duplicate it. The copy lands right after it and keeps the OLD
#     "总计" data/format; this copy will become the NEW "总计" sheet (sheetId 6).
$total.Copy($null, $total)
$newTotal = $wb.Worksheets.Item($total.Index + 1)

# --- Step 2: rename sheets. Rename the original first so the name "总计" is
#     free for the copy to take.
$total.Name = "2022-Q1"
$newTotal.Name = "总计"

# ======================================================================
# Build the new "2022-Q1" sheet: fund-holding detail (same shape as the
# other quarterly sheets: B:H headers, rows 2-4 data, A2:A4 index column).
# Use the "2021-Q2" sheet (3rd tab) as a formatting template for A1:H4,
# then overwrite every value with the 2022-Q1 figures.
# ======================================================================
$fmtSrc = $wb.Worksheets.Item(3)
$total.Cells.Clear()
$fmtSrc.Range("B1:H4").Copy($total.Range("B1:H4"))
$fmtSrc.Range("A2:A4").Copy($total.Range("A2:A4"))

# Header row
$total.Cells.Item(1,2).Value = "基金代码"
$total.Cells.Item(1,3).Value = "基金名称"
$total.Cells.Item(1,4).Value = "基金规模"
$total.Cells.Item(1,5).Value = "股票总仓位"
$total.Cells.Item(1,6).Value = "仓位占比"
$total.Cells.Item(1,7).Value = "持有市值(亿元)"
$total.Cells.Item(1,8).Value = "仓位排名"

# Helper to write a text value that must stay text even though it looks
# like a number (fund code / size / weight figures are stored as text in
# the source workbook).
function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 2 - 160211 国泰中小盘成长混合(LOF)
Set-TextCell $total.Cells.Item(2,2) "160211"
Set-TextCell $total.Cells.Item(2,3) "国泰中小盘成长混合(LOF)"
Set-TextCell $total.Cells.Item(2,4) "6.78"
Set-TextCell $total.Cells.Item(2,5) "89.07"
Set-TextCell $total.Cells.Item(2,6) "2.74"
Set-TextCell $total.Cells.Item(2,7) "0.1858"
$total.Cells.Item(2,8).Value = 10

# Row 3 - 006267 诺德量化核心灵活配置混合A
Set-TextCell $total.Cells.Item(3,2) "006267"
Set-TextCell $total.Cells.Item(3,3) "诺德量化核心灵活配置混合A"
Set-TextCell $total.Cells.Item(3,4) "1.84"
Set-TextCell $total.Cells.Item(3,5) "93.91"
Set-TextCell $total.Cells.Item(3,6) "2.57"
Set-TextCell $total.Cells.Item(3,7) "0.0473"
$total.Cells.Item(3,8).Value = 10

# Row 4 - 006268 诺德量化核心灵活配置混合C
Set-TextCell $total.Cells.Item(4,2) "006268"
Set-TextCell $total.Cells.Item(4,3) "诺德量化核心灵活配置混合C"
Set-TextCell $total.Cells.Item(4,4) "0.50"
Set-TextCell $total.Cells.Item(4,5) "93.91"
Set-TextCell $total.Cells.Item(4,6) "2.57"
Set-TextCell $total.Cells.Item(4,7) "0.0128"
$total.Cells.Item(4,8).Value = 10

# ======================================================================
# Update the (new, renamed-copy) "总计" sheet: prepend a 2022-Q1 summary
# row above the existing quarters, shifting everything else down by one
# row and re-indexing column A (0,1,2,3,4).
# ======================================================================

# Capture the old data (rows 2..5, columns B..D) before it gets shifted.
$oldB = @()
$oldC = @()
$oldD = @()
for ($i = 2; $i -le 5; $i++) {
    $oldB += $newTotal.Cells.Item($i,2).Value2
    $oldC += $newTotal.Cells.Item($i,3).Value2
    $oldD += $newTotal.Cells.Item($i,4).Value2
}

# Extend the A-column index style down to the new last row (row 6).
$newTotal.Range("A2:A5").Copy($newTotal.Range("A3:A6"))

# Re-write rows 3..6 with the old rows 2..5 (shifted down by one), and
# re-number the A-column index 0,1,2,3,4 top to bottom.
for ($i = 0; $i -le 3; $i++) {
    $r = 6 - $i
    $newTotal.Cells.Item($r,1).Value = $r - 2
    $newTotal.Cells.Item($r,2).Value = $oldB[3-$i]
    $newTotal.Cells.Item($r,3).Value = $oldC[3-$i]
    $newTotal.Cells.Item($r,4).Value = $oldD[3-$i]
}

# New row 2: 2022-Q1 summary.
$newTotal.Cells.Item(2,1).Value = 0
$newTotal.Cells.Item(2,2).Value = "2022-Q1"
$newTotal.Cells.Item(2,3).Value = 3
$newTotal.Cells.Item(2,4).Value = 0.25
